# Deleted Drivers Script added
#
# - Configuration!B3 ("Send Mail") flips from "no" to "yes".
# - The previously-active sheet (Company_Profile) is left with its
#   selection moved to E7 and is no longer the active tab.
# - The Configuration sheet becomes the active tab, with C4 selected.

$wb = $excel.ActiveWorkbook

$wsConfig  = $wb.Worksheets.Item("Configuration")
$wsCompany = $wb.Worksheets.Item("Company_Profile")

# Flip the "Send Mail" flag.
$wsConfig.Range("B3").Value = "yes"

# Update the selection left behind on Company_Profile before switching away.
$wsCompany.Activate()
$wsCompany.Range("E7").Select()

# Switch to Configuration and leave C4 selected as the new active cell.
$wsConfig.Activate()
$wsConfig.Range("C4").Select()
